# 质检费导入模板.xlsx -- remove the "商家编号" column and turn the
# "商家名称" header into a required "*商家名称" field (red, centered).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the header row one column to the left starting at B1, which
# overwrites "商家编号" with the contents that used to sit in the next
# column over (and so on down the row). This mirrors deleting column B
# without touching the <cols> width/style definitions.
$ws.Range("B1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = $ws.Range("E1").Value()
$ws.Range("E1").Value = $ws.Range("F1").Value()
$ws.Range("F1").Value = $ws.Range("G1").Value()
$ws.Range("G1").Value = $ws.Range("H1").Value()
$ws.Range("H1").Value = $ws.Range("I1").Value()
$ws.Range("I1").Value = $ws.Range("J1").Value()

# The last header ("所属机构") lands on I1 with a plain/default look
# (matches the source workbook, where column J had no explicit style).
$blank = $ws.Cells.Item(5, 5).Style()
$ws.Range("I1").Style = $blank

# J1 is now unused entirely -- drop it instead of leaving an empty cell.
$ws.Range("J1").Clear()

# Rename the merchant-name header to flag it as a required field and
# give it its own styling: red 宋体 text, centered.
$ws.Range("B1").Value = "*商家名称"
$ws.Range("B1").Font.Color = 255
$ws.Range("B1").HorizontalAlignment = -4108

# Restore the selection Excel left the workbook with after the edit.
$ws.Range("D10").Select()
